$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.817.30"
$ws.Range("E2").Value = "  +0.88%  "

$ws.Range("D3").Value = "3.954.91"
$ws.Range("E3").Value = "  -2.22%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "610.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +13.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +11.22%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.677"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.89%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.748"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.50%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.180"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.60%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "56.10"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.83%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000332"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.31%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.02"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.42%  "

$ws.Range("D14").Value = "4.596.16"
$ws.Range("E14").Value = "  -1.95%  "

$ws.Range("D15").Value = "3.976.67"

$ws.Range("E16").Value = "  +2.51%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.98"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.54%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.38"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.02%  "

$ws.Range("E19").Value = "  -0.40%  "

$ws.Range("D20").Value = "72.802.64"
$ws.Range("E20").Value = "  +0.92%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "435.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.67%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.85"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +13.80%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "95.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.46%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.35"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.95%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.72%  "

$ws.Range("E26").Value = "  -8.97%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.41%  "

$ws.Range("E28").Value = "  -0.03%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.38"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.05%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.82"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.81%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.91"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.77%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.60"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.73%  "

$ws.Range("E33").Value = "  -3.62%  "

$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "47.65"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.54%  "

$ws.Range("B35").Value = "PEPE"
$ws.Range("C35").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0000101"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +10.44%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "70.09"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.62%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "638.88"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.58%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.429"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.91%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.45"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.75%  "

$ws.Range("E40").Value = "  -1.13%  "

$ws.Range("E41").Value = "  -0.07%  "

$ws.Range("E42").Value = "  +0.23%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.78"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.65%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0481"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.12%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.19"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.09%  "

$ws.Range("E46").Value = "  -1.83%  "

$ws.Range("E47").Value = "  +0.56%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.59"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.41%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.85"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +25.32%  "

$ws.Range("D50").Value = "2.852.20"
$ws.Range("E50").Value = "  +3.31%  "

$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.54%  "
